$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric "time[us]" column updates (column S)
$ws.Range("S6").Value = 1664
$ws.Range("S114").Value = 186810

# Layer/tensor-name column updates (column B).
# The new values are numeric-looking ONNX tensor ids (e.g. "701") that must
# remain stored as text, matching the original inline-string cell type.
# Writing them through a text-valued formula + Copy/PasteSpecial(values)
# keeps them as text without Excel re-interpreting them as numbers.
$helper = $ws.Range("Z1")
$map = [ordered]@{
    "B8" = "701"
    "B9" = "704"
    "B10" = "707"
    "B11" = "710"
    "B12" = "457"
    "B13" = "713"
    "B14" = "716"
    "B15" = "719"
    "B16" = "467"
    "B17" = "722"
    "B18" = "725"
    "B19" = "474"
    "B20" = "728"
    "B21" = "731"
    "B22" = "734"
    "B23" = "484"
    "B24" = "737"
    "B25" = "740"
    "B26" = "491"
    "B27" = "743"
    "B28" = "746"
    "B29" = "498"
    "B30" = "749"
    "B31" = "752"
    "B32" = "505"
    "B33" = "755"
    "B34" = "758"
    "B35" = "512"
    "B36" = "761"
    "B37" = "764"
    "B38" = "519"
    "B39" = "767"
    "B40" = "770"
    "B41" = "526"
    "B42" = "773"
    "B43" = "776"
    "B44" = "533"
    "B45" = "779"
    "B46" = "782"
    "B47" = "785"
    "B48" = "543"
    "B49" = "788"
    "B50" = "791"
    "B51" = "550"
    "B52" = "794"
    "B53" = "797"
    "B54" = "557"
    "B55" = "800"
    "B56" = "803"
    "B57" = "564"
    "B58" = "806"
    "B59" = "809"
    "B60" = "571"
    "B61" = "812"
    "B62" = "815"
    "B63" = "578"
    "B64" = "818"
    "B65" = "821"
    "B66" = "585"
    "B67" = "824"
    "B68" = "827"
    "B69" = "592"
    "B70" = "830"
    "B71" = "833"
    "B72" = "836"
    "B73" = "602"
    "B74" = "839"
    "B75" = "842"
    "B76" = "609"
    "B77" = "845"
    "B78" = "848"
    "B79" = "616"
    "B80" = "851"
    "B81" = "854"
    "B82" = "623"
    "B83" = "857"
    "B84" = "860"
    "B85" = "863"
    "B86" = "866"
    "B87" = "869"
    "B88" = "872"
    "B90" = "875"
    "B91" = "651"
    "B92" = "878"
    "B93" = "881"
    "B94" = "884"
    "B95" = "887"
    "B96" = "890"
    "B97" = "893"
    "B99" = "896"
    "B100" = "680"
    "B101" = "899"
    "B102" = "902"
    "B103" = "905"
    "B104" = "908"
    "B105" = "911"
    "B106" = "914"
}
foreach ($addr in $map.Keys) {
    $helper.Formula = '="' + $map[$addr] + '"'
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}
$helper.Clear()
$excel.CutCopyMode = $false
